$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)

# Resize / reposition the subtitle placeholder
$sh.Left = 48.15433070866142
$sh.Width = 657.7108661417323

$tr = $sh.TextFrame.TextRange
$tr.Text = "Chương 14. Phân lớp và ứng dụng trong tìm "
$tr.InsertAfter("kiếm") | Out-Null
$tr.InsertAfter("`rIIR.C13. Text classification and Naive Bayes") | Out-Null

$tr.Paragraphs(1).ParagraphFormat.Alignment = 4
$tr.Paragraphs(2).ParagraphFormat.Alignment = 4
